# StructureDefinition-SegundoApellido.xlsx — "version final sin errores"
#
# 1. Metadata sheet: bump Version 0.4.0 -> 0.7.0, refresh the publish Date,
#    and correct the extension Context from "element:Element" to
#    "element:Patient".
# 2. Elements sheet: the root "Extension" row (row 1) was missing its
#    invariants; add the same ele-1 / ext-1 FHIR invariant text that already
#    appears on the "Element.extension" row.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.7.0"
$meta.Range("B8").Value = "2023-09-13T17:11:14-03:00"
$meta.Range("B20").Value = "element:Patient"

$elements = $wb.Worksheets.Item("Elements")
$invariants = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + "`n" + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$elements.Range("AJ1").Value = $invariants
